# Resolve the active workbook/worksheet handed to us by the COM host.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 36's height ---
# It was stuck at 45pt despite its wrapped abstract text needing much more
# room (consistent with the autofit height used by every other data row);
# correct it to 315pt to match the rest of the sheet's sizing convention.
$ws.Rows.Item(36).RowHeight = 315

# --- Append 10 new article records (rows 42-51) ---
# Columns: A=pii, B=title, C=abstract, D=keywords
# (values are entered title-first, then pii/abstract/keywords, to match
# the order the new entries were typed/pasted into the sheet)

# Row 42: S0306457319303474
$ws.Cells.Item(42, 2).Value = 'Network measures: A new paradigm towards reliable novel word sense detection'
$ws.Cells.Item(42, 1).Value = 'S0306457319303474'
$ws.Cells.Item(42, 3).Value = 'In this era of digitization, with the fast flow of information on the web, words are being used to denote newer meanings. Thus novel sense detection becomes a crucial and challenging task in order to build any natural language processing application which depends on the efficient semantic representation of words. With the recent availability of large amounts of digitized texts, automated analysis of language evolution has become possible. Given corpus from two different time periods, the main focus of our work is to detect the words evolved with a novel sense precisely. We pose this problem as a binary classification task to detect whether a new sense of a target word has emerged. This paper presents a unique proposal based on network features to improve the precision of this task of detecting emerged new sense of a target word. For a candidate word where a new sense has been detected by comparing the sense clusters induced at two different time periods, we further compare the network properties of the subgraphs induced from novel sense clusters across these two time periods. Using the mean fractional change in edge density, structural similarity and average path length as features in a Support Vector Machine (SVM) classifier, manual evaluation gives precision values of 0.86 and 0.74 for the task of new sense detection, when tested on 2 distinct time-point pairs, in comparison to the precision values in the range of 0.23-0.32, when the proposed scheme is not used. The outlined method can, therefore, be used as a new post-hoc step to improve the precision of novel word sense detection in a robust and reliable way where the underlying framework uses a graph structure. Another important observation is that even though our proposal is a post-hoc step, it can be used in isolation and that itself results in a very decent performance achieving a precision of 0.54-0.62. Finally, we also show that our method is able to detect well-known historical shifts in 80% cases.'
$ws.Cells.Item(42, 4).Value = 'Novel sense detection; Distributional thesaurus network; Complex network measures'
$ws.Rows.Item(42).RowHeight = 345

# Row 43: S0950705118301539
$ws.Cells.Item(43, 2).Value = 'Multilevel approach for combinatorial optimization in bipartite network'
$ws.Cells.Item(43, 1).Value = 'S0950705118301539'
$ws.Cells.Item(43, 3).Value = 'Multilevel approaches aim at reducing the cost of a target algorithm over a given network by applying it to a coarsened (or reduced) version of the original network. They have been successfully employed in a variety of problems, most notably community detection. However, current solutions are not directly applicable to bipartite networks and the literature lacks studies that illustrate their application for solving multilevel optimization problems in such networks. This article addresses this gap and introduces a multilevel optimization approach for bipartite networks and the implementation of a general multilevel framework including novel algorithms for coarsening and uncorsening, applicable to a variety of problems. We analyze how the proposed multilevel strategy affects the topological features of bipartite networks and show that a controlled coarsening strategy can preserve properties such as degree and clustering coefficient centralities. The applicability of the general framework is illustrated in two optimization problems, one for solving the Barber’s modularity for community detection and the second for dimensionality reduction in text classification. We show that the solutions thus obtained are statistically equivalent, regarding accuracy, to those of conventional approaches, whilst requiring considerably lower execution times.'
$ws.Cells.Item(43, 4).Value = 'Complex networks; Bipartite networks; Combinatorial optimization; Meta-heuristic; Multilevel optimization; Large-scale networks'
$ws.Rows.Item(43).RowHeight = 240

# Row 44: S0378437113005839
$ws.Cells.Item(44, 2).Value = 'Lobby index as a network centrality measure'
$ws.Cells.Item(44, 1).Value = 'S0378437113005839'
$ws.Cells.Item(44, 3).Value = 'We study the lobby index (l-index for short) as a local node centrality measure for complex networks. The l-index is compared with degree (a local measure), betweenness and Eigenvector centralities (two global measures) in the case of a biological network (Yeast interaction protein–protein network) and a linguistic network (Moby Thesaurus II). In both networks, the l-index has a poor correlation with betweenness but correlates with degree and Eigenvector centralities. Although being local, the l-index carries more information about its neighbors than degree centrality. Also, it requires much less time to compute when compared with Eigenvector centrality. Results show that the l-index produces better results than degree and Eigenvector centrality for ranking purposes.'
$ws.Cells.Item(44, 4).Value = 'Lobby index; Centrality; Degree; Betweenness; Eigenvector; Hirsch index'
$ws.Rows.Item(44).RowHeight = 135

# Row 45: S0957417421002748
$ws.Cells.Item(45, 2).Value = 'A framework for inventor collaboration recommendation system based on network approach'
$ws.Cells.Item(45, 1).Value = 'S0957417421002748'
$ws.Cells.Item(45, 3).Value = 'Precise and timely information about opportunities for potential collaborations is very vital for the collaboration-intense research environment prevailing in innovation ecosystems. As the identification of suitable inventors for collaboration will be decisive for inventors in different phases of their careers, inventor collaboration recommendation systems are of great importance. Existing recommendation system frameworks for collaboration recommendations for academic authors and inventors are slightly intensive on the usage of link semantics. Like academic collaboration through co-authorship, collaborations of inventors through co-inventorship of patents can be found in almost all industrial areas in various degrees. Network representation of co-inventorship can be used to retrieve many insights that can even be vital for policymaking. In this work, for inventor collaboration recommendations, a minimal link semantics (MLS) approach based framework is built to overcome these major drawbacks and to improve usability. The case of inventors in the area ‘Wireless power transmission’ is analyzed using patent data for the demonstration of the MLS framework and on evaluation, the framework is found to be capable of retrieving novel and diverse recommendations to and from inventors that belong to different phases of a career.'
$ws.Cells.Item(45, 4).Value = 'Complex networks; Patent analysis; Inventor Collaboration; Co-inventorship; Link prediction; Recommendation system'
$ws.Rows.Item(45).RowHeight = 240

# Row 46: S0378437116305234
$ws.Cells.Item(46, 2).Value = 'The rapid bi-level exploration on the evolution of regional solar energy development'
$ws.Cells.Item(46, 1).Value = 'S0378437116305234'
$ws.Cells.Item(46, 3).Value = 'As one of the renewable energy, solar energy is experiencing increased but exploratory development worldwide. The positive or negative influences of regional characteristics, like economy, production capacity and allowance policies, make them have uneven solar energy development. In this paper, we aim at quickly exploring the features of provincial solar energy development, and their concerns about solar energy. We take China as a typical case, and combine text mining and two-actor networks. We find that the classification of levels based on certain nodes and the amount of degree avoids missing meaningful information that may be ignored by global level results. Moreover, eastern provinces are hot focus for the media, western countries are key to bridge the networks and special administrative region has local development features; third, most focus points are more about the application than the improvement of material. The exploration of news provides practical information to adjust researches and development strategies of solar energy. Moreover, the bi-level exploration, which can also be expanded to multi-level, is helpful for governments or researchers to grasp more targeted and precise knowledge.'
$ws.Cells.Item(46, 4).Value = 'Solar energy; Regional development; Text mining; Complex network'
$ws.Rows.Item(46).RowHeight = 210

# Row 47: S0020025514004137
$ws.Cells.Item(47, 2).Value = 'Linguistic performance evaluation for an ERP system with link failures'
$ws.Cells.Item(47, 1).Value = 'S0020025514004137'
$ws.Cells.Item(47, 3).Value = 'An Enterprise Resource Planning (ERP) system is a complex network composed of various business processes. It can be called an ERP net. This paper proposes an analytic method to evaluate the Linguistic performance of such net under link failure situations. A link failure in an ERP net means that the software or hardware between processes may malfunction. To facility such evaluation, the nodes in the net denote the persons responsible for the business tasks during the processes. The links between nodes denote the process precedence relationships in the ERP system. When the process starts, the documents (jobs) are initiated from the source node to its succeeding nodes. Finally, the documents are released in the destination node. Thus, the performance of an ERP system is related to the document flow under the net. The performance failure of an ERP system is therefore defined by the condition that the document flow of the system is under the acceptable level d. By using the fuzzy linguistic results of the ERP examination of the users, we propose a fuzzy linguistic performance index, defuzzified from the probability of maximal flow not less than d, to evaluate the performance of an ERP system. An algorithm is subsequently proposed to generate the performance index under link failure situations, which can be used to real time assess the system performance either before or after the system going live.'
$ws.Cells.Item(47, 4).Value = 'Enterprise resource planning; ERP net; Fuzzy mathematics; Performance evaluation; Minimal path'
$ws.Rows.Item(47).RowHeight = 255

# Row 48: S095741742200094X
$ws.Cells.Item(48, 2).Value = 'A network-based feature extraction model for imbalanced text data'
$ws.Cells.Item(48, 1).Value = 'S095741742200094X'
$ws.Cells.Item(48, 3).Value = 'The explosive growth of text data has attracted many researchers to explore the efficient method to extract valuable hidden information. Many technologies, especially deep learning methods, have achieved great success in text analysis. However, the most powerful methods always require a considerable quantity of data for training, which may suffer from imbalanced data in some cases. In this paper, we propose a network-based Convolution Neural Network (NCNN) to mitigate the effect of imbalanced data. The proposed model first generates new synthetic samples for the imbalanced data based on the random walking of the network. Then an extra layer called Polar Layer is introduced to connect the output from the network model of the text to the classical CNN. Two electing strategies (n-NCNN and x-NCNN) are proposed to improve the performance of NCNN further. In the experimental section, the proposed model is applied to Reuters 21578 and WebKb. By comparing with six approaches, we prove the effectiveness of the proposed NCNN model on the imbalanced text data.'
$ws.Cells.Item(48, 4).Value = 'Complex Network; CNN; Text Analysis; Imbalanced Data; Random Walk'
$ws.Rows.Item(48).RowHeight = 195

# Row 49: S0378437113001088
$ws.Cells.Item(49, 2).Value = 'A network approach based on cliques'
$ws.Cells.Item(49, 1).Value = 'S0378437113001088'
$ws.Cells.Item(49, 3).Value = 'The characterization of complex networks is a procedure that is currently found in several research studies. Nevertheless, few studies present a discussion on networks in which the basic element is a clique. In this paper, we propose an approach based on a network of cliques. This approach consists not only of a set of new indices to capture the properties of a network of cliques but also of a method to characterize complex networks of cliques (i.e., some of the parameters are proposed to characterize the small-world phenomenon in networks of cliques). The results obtained are consistent with results from classical methods used to characterize complex networks.'
$ws.Cells.Item(49, 4).Value = 'Networks of cliques; Complex networks; Small-world phenomenon; Social network; Semantic networks'
$ws.Rows.Item(49).RowHeight = 120

# Row 50: S0378437118310434
$ws.Cells.Item(50, 2).Value = 'A framework of community detection based on individual labels in attribute networks'
$ws.Cells.Item(50, 1).Value = 'S0378437118310434'
$ws.Cells.Item(50, 3).Value = 'Community detection is an important problem for understanding the structure and function of complex networks and has attracted a lot of attention in recent decades. Most community detection algorithms only focus on the topology of networks. However, there is still much valuable information hidden in the networks, such as the attributes or content of the nodes and the useful prior information. Obviously, taking full advantage of these resources can improve the effectiveness of community detection. In this paper, we present a semi-supervised community detection framework named SCDAN (Semi-supervised Community Detection in Attribute Networks), in which a non-negative matrix factorization model is utilized to effectively integrate network topology, node attributes and individual labels simultaneously. The comparative experiments on real-world networks show that SCDAN significantly improves the performance of community detection and provides semantic interpretation of communities.'
$ws.Cells.Item(50, 4).Value = 'Community detection; Attribute network; Individual label; Non-negative matrix factorization'
$ws.Rows.Item(50).RowHeight = 180

# Row 51: S095032931831022X
$ws.Cells.Item(51, 2).Value = 'Consumers’ associative networks of plant-based food product communications'
$ws.Cells.Item(51, 1).Value = 'S095032931831022X'
$ws.Cells.Item(51, 3).Value = 'Food producers respond to the current consumer trend of clean label products and reducing meat consumption by increasingly offering plant-based food products and transparent, understandable ingredient lists. However, consumer interest can be driven by various motives and food producers face the challenge of identifying the most effective motive to address. We analyze concept maps of 90 consumers who received information that positioned plant-based food products as sustainable, healthy, or with a transparent ingredient focus. To assess the applicability of text mining with a view to reducing coder bias and the duration of qualitative data analysis, we compared the results of text mining versus a human coder approach. Our results show that human coder analysis results in more detail, however the advantage of the text mining procedure is that it can run independently and analyze qualitative data more objectively. When a high degree of control and depth of analysis is necessary to satisfy the study objective, human coding might have its rewards. For the current study, both approaches draw a similar picture of the associative networks and are therefore equally suitable to satisfy the study objective. When plant-based diets are communicated solely based on the ingredient used for substituting animal-based ingredients, associative networks are less complex and associations are primarily concerned with taste. A health communication perspective results in more complex networks with a focus on other food product properties such as processing degree and nutrition. A sustainability communication also results in higher complexity, with fewer associations concerning the product properties itself, but rather with the environmental impact and the authenticity of the product. The in-depth understanding of consumers’ associations evoked by communicating different perspectives of plant-based food products can be used by practitioners in tailoring their marketing activities to the characteristics of their product offerings.'
$ws.Cells.Item(51, 4).Value = 'Concept mapping; Flexitarians; Text mining; Network analysis; Potato protein'
$ws.Rows.Item(51).RowHeight = 360

# --- Match the workbook's final on-screen selection/scroll state ---
# After appending the new rows, the author's cursor ended on the last
# populated cell (D51), with the frozen header pane scrolled so row 51
# is the first visible data row.
$ws.Range("D51").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 51
$excel.ActiveWindow.ScrollColumn = 1

